# Reorder the "Recorded By" (column G) comma-separated list so that any
# email-like entries (containing "@") come first, followed by the
# remaining entries (e.g. "system"/"System"), each group preserving its
# original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row()
$lastRow = $firstRow + $used.Rows.Count() - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value()

    if ($val -eq $null) {
        continue
    }

    $text = [string]$val
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ','
    $emails = @()
    $others = @()

    foreach ($part in $parts) {
        $trimmed = $part.Trim()
        if ($trimmed -like "*@*") {
            $emails += $trimmed
        } else {
            $others += $trimmed
        }
    }

    $newParts = $emails + $others
    $newText = [string]::Join(", ", $newParts)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
